$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 values to "custom accuracy" (rounded to 2 decimal places) ---
$row5 = @{
    2  = 14.27                 # B5
    3  = 10.69                 # C5
    4  = 0.9                   # D5
    5  = 31.29                 # E5
    6  = 25.45                 # F5
    7  = 10.9                  # G5
    8  = 45.11                 # H5
    9  = 17.32                 # I5
    10 = 8.07                  # J5
    11 = 11.18                 # K5
    12 = 12.57                 # L5
    13 = 13.43                 # M5
    14 = 3.81                  # N5
    15 = 11.28                 # O5
    16 = 15.86                 # P5
    17 = 9.619999999999999     # Q5
    18 = 0.26                  # R5
    19 = 0.6                   # S5
    20 = 164.36                # T5
    21 = 31.58                 # U5
    22 = 10.41                 # V5
    23 = 21.09                 # W5
    24 = 11.25                 # X5
    25 = 1.48                  # Y5
    26 = 22.38                 # Z5
    27 = 9.140000000000001     # AA5
    28 = 8.199999999999999     # AB5
    29 = 9.699999999999999     # AC5
    30 = 13.27                 # AD5
    31 = 0.47                  # AE5
    32 = 41.16                 # AF5
    33 = 5.9                   # AG5
    34 = 12.96                 # AH5
}

foreach ($col in $row5.Keys) {
    $ws.Cells.Item(5, $col).Value = $row5[$col]
}

# --- Remove the last data row (row 6), shrinking the used range to A1:AH5 ---
$ws.Rows("6").Delete()

# --- Narrow column V (22nd column) from width 8 to width 7 ---
$ws.Columns("V").ColumnWidth = 6.17
